# Apply the "updated with tests, new code (lstm)" edit:
#  - Rename several header strings (column labels) in row 1
#  - Replace column C ("GDP", formerly "Gross_National_Income") values
#    with new model-predicted figures for every data row (2-53)
#  - Row 47's column C was a placeholder text value ".."; it now gets a
#    real numeric GDP estimate like every other row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) values, rows 2-53 ---
$cValues = @{
    2  = 2771.04675450926
    3  = 5596.139681459835
    4  = 2934.187009790061
    5  = 2870.311589353206
    6  = 1460.056109840828
    7  = 10594.98659239237
    8  = 4547.50930098406
    9  = 4729.735976516416
    10 = 5730.354774594881
    11 = 2983.242707849043
    12 = 2898.942214704482
    13 = 1503.870423231357
    14 = 4633.590358399045
    15 = 5082.354756663512
    16 = 2948.84548976845
    17 = 5885.254624554112
    18 = 3083.80337578809
    19 = 2965.153206179127
    20 = 1577.487171555845
    21 = 5360.226632400601
    22 = 4921.848409120176
    23 = 6051.685746144485
    24 = 3156.723844635973
    25 = 2999.422762626143
    26 = 1657.651524528445
    27 = 5122.180090208862
    28 = 2995.45235738661
    29 = 6203.843262938323
    30 = 3212.740625904757
    31 = 3056.152683606517
    32 = 1716.389195271215
    33 = 5295.682695961288
    34 = 3087.12349650562
    35 = 6255.426161047989
    36 = 5412.131646018807
    37 = 3252.634165082374
    38 = 3137.260298393558
    39 = 1775.027517189621
    40 = 5996.49696468919
    41 = 6522.736799041846
    42 = 5330.539154475424
    43 = 3314.741082534716
    44 = 3210.869677115934
    45 = 1836.014008604312
    46 = 6114.227214287786
    47 = 6550.274372976741
    48 = 5176.058803160127
    49 = 3382.563653843273
    50 = 3242.636921959078
    51 = 3212.81539531051
    52 = 1895.214690888655
    53 = 6262.368904654469
}

foreach ($row in $cValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
}
